# Update the 南宁-漫展信息 workbook: both the "展览" and "全部类型" sheets list
# upcoming Nanning comic-expo events. The event that already happened
# (2024-03-30, cancelled ANE/DACG show) is dropped, every remaining event
# shifts up one row, and a new event (2024-06-09, 布谷鸟动漫展4th) is
# appended at the end -- with a couple of the "想去人数" (F) figures
# refreshed to newer counts.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 2
    $ws.Range("B2").Value = "'2024-04-11"
    $ws.Range("C2").Value = "南宁·三月三漫次元国风动漫节"
    $ws.Range("D2").Value = "亭洪路45号 百益上河城"
    $ws.Range("E2").Value = "2024.04.11 10:00-04.12 17:00"
    $ws.Range("F2").Value = 448
    $ws.Range("G2").Value = 45
    $ws.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=83139"
    $ws.Range("I2").Value = "//i2.hdslb.com/bfs/openplatform/202403/nqZxtIdd1711001896521.jpeg"

    # Row 3
    $ws.Range("B3").Value = "'2024-05-01"
    $ws.Range("C3").Value = "南宁·2024三月三国潮动漫节（良牙春典）"
    $ws.Range("D3").Value = "民族大道106号 南宁国际会展中心"
    $ws.Range("E3").Value = "2024.05.01 09:30-05.02 17:30"
    $ws.Range("F3").Value = 3243
    $ws.Range("G3").Value = 55
    $ws.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=82416"
    $ws.Range("I3").Value = "//i0.hdslb.com/bfs/openplatform/202403/b3YxmMm81711075370604.jpeg"

    # Row 4
    $ws.Range("B4").Value = "'2024-05-19"
    $ws.Range("C4").Value = "南宁·原x穹x崩only"
    $ws.Range("D4").Value = "明秀东路157号 利泰国际大酒店"
    $ws.Range("E4").Value = "2024.05.19 10:00-05.19 17:00"
    $ws.Range("F4").Value = 79
    $ws.Range("G4").Value = 35
    $ws.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=83070"
    $ws.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202403/I8tScigE1710918412731.jpeg"

    # Row 5
    $ws.Range("B5").Value = "'2024-06-09"
    $ws.Range("C5").Value = "南宁·布谷鸟动漫展4th"
    $ws.Range("D5").Value = "亭洪路45号 百益上河城"
    $ws.Range("E5").Value = "2024.06.09 10:00-06.10 17:00"
    $ws.Range("F5").Value = 647
    $ws.Range("G5").Value = 35
    $ws.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=82241"
    $ws.Range("I5").Value = "//i2.hdslb.com/bfs/openplatform/202403/uzZqZov91709281147333.jpeg"

    # The old row 6 (布谷鸟动漫展4th, now merged into row 5 above) is removed
    # entirely so the sheet shrinks from A1:I6 back down to A1:I5.
    $ws.Rows.Item(6).Delete()
}
